$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO" (columns A:R) ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Insert a new row at 248, shifting the existing rows 248-275 down to 249-276.
$ws1.Rows.Item(248).Insert()

# Populate the newly inserted row with the new asesor/cliente pair.
$ws1.Cells.Item(248, 1).Value = "OFICINA-CATAECSA"
$ws1.Cells.Item(248, 2).Value = "MOROCHO PLAZA SHIRLEY AURELIA"
for ($c = 3; $c -le 18; $c++) {
    $ws1.Cells.Item(248, $c).Value = 0
}

# The final summary row (now row 276) still refers to the old total count of
# 273 clients in its "N de 273" labels; bump that to 274 everywhere.
for ($c = 3; $c -le 18; $c++) {
    $cell = $ws1.Cells.Item(276, $c)
    $old = $cell.Value()
    $cell.Value = ($old -replace "273", "274")
}

# --- Sheet "VENTA MENSUAL" (columns A:G) ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Same insertion at row 248.
$ws2.Rows.Item(248).Insert()

$ws2.Cells.Item(248, 1).Value = "OFICINA-CATAECSA"
$ws2.Cells.Item(248, 2).Value = "MOROCHO PLAZA SHIRLEY AURELIA"
for ($c = 3; $c -le 7; $c++) {
    $ws2.Cells.Item(248, $c).Value = 0
}

# The totals row on this sheet (now row 276) holds numeric sums only, no
# "de N" labels, so its values are unaffected by the insert and need no edit.
